$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.925.24"
$ws.Range("E2").Value = "  -2.36%  "

$ws.Range("D3").Value = "3.408.35"
$ws.Range("E3").Value = "  -3.00%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("E5").Value = "  -2.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.91"
$ws.Range("E6").Value = "  -5.68%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "3.409.16"
$ws.Range("E8").Value = "  -2.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("E9").Value = "  -2.20%  "

$ws.Range("E10").Value = "  -1.94%  "

$ws.Range("E11").Value = "  -2.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.381"
$ws.Range("E12").Value = "  -1.42%  "

$ws.Range("D13").Value = "3.991.28"
$ws.Range("E13").Value = "  -2.87%  "

$ws.Range("E14").Value = "  -0.82%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000174"
$ws.Range("E15").Value = "  -3.89%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.408.54"
$ws.Range("E16").Value = "  -3.00%  "

$ws.Range("D17").Value = "62.981.63"
$ws.Range("E17").Value = "  -2.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.05"
$ws.Range("E18").Value = "  -3.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.60"
$ws.Range("E19").Value = "  -3.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.72"
$ws.Range("E20").Value = "  -0.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.22"
$ws.Range("E21").Value = "  -2.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "378.35"
$ws.Range("E22").Value = "  -4.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.561"
$ws.Range("E23").Value = "  -2.67%  "

$ws.Range("D24").Value = "3.545.93"
$ws.Range("E24").Value = "  -2.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.59"
$ws.Range("E25").Value = "  -2.93%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("E27").Value = "  -7.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("E29").Value = "  -5.33%  "

$ws.Range("E30").Value = "  -4.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.92"
$ws.Range("E31").Value = "  -4.57%  "

$ws.Range("E32").Value = "  -3.99%  "

$ws.Range("E33").Value = "  -2.53%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").Value = "3.435.16"
$ws.Range("E35").Value = "  -2.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.89"
$ws.Range("E36").Value = "  -2.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.33"
$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.78"
$ws.Range("E38").Value = "  -2.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "163.94"
$ws.Range("E39").Value = "  -2.06%  "

$ws.Range("E40").Value = "  -3.26%  "

$ws.Range("E41").Value = "  -3.67%  "

$ws.Range("E42").Value = "  -3.58%  "

$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.65"
$ws.Range("E44").Value = "  -2.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.31"
$ws.Range("E45").Value = "  -3.30%  "

$ws.Range("E46").Value = "  -4.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.91"
$ws.Range("E47").Value = "  -9.99%  "

$ws.Range("E48").Value = "  -6.88%  "

$ws.Range("E49").Value = "  -1.52%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.260.99"
$ws.Range("E50").Value = "  -5.69%  "

$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.862"
$ws.Range("E51").Value = "  -4.09%  "
